$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Register / Add-user testcases review rows (72-79): reviewer
#    "Tarek" finished reviewing these rows -> fill in the Reviewer
#    (H) and Status (I) columns that were previously left blank.
#    H gets the sheet/row default formatting (typed fresh), I is
#    copied from the existing "Done" status cell at I70 so it picks
#    up the same style already used for a "Done" status in this
#    document.
# ------------------------------------------------------------------
$ws.Range("I70").Copy() | Out-Null
$ws.Range("I72:I79").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

for ($r = 72; $r -le 79; $r++) {
    $ws.Cells.Item($r, 8).Value = "Tarek"   # H: Reviewer
    $ws.Cells.Item($r, 9).Value = "Done"    # I: Status
}

# ------------------------------------------------------------------
# 2) Review history is updated: four new review-comment rows (92-95)
#    added at the bottom of the log, following the same layout/
#    styling as the row directly above them (row 91), with a
#    slightly smaller default font that the new rows were entered
#    with. The highlighted "V1" badge style (already used on B91,
#    D91 and E91) is carried across the whole new row, not just the
#    badge/date columns.
# ------------------------------------------------------------------
$ws.Range("B91").Copy() | Out-Null
$ws.Range("A92:A95").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C92:C95").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F92:F95").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("G92:G95").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B92:B95").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D91").Copy() | Out-Null
$ws.Range("D92:D95").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E91").Copy() | Out-Null
$ws.Range("E92:E95").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# New rows use a smaller (default-ish) font than row 91 for every
# column except the comment body (D), which uses a 12pt font.
$ws.Range("A92:C95").Font.Size = 11
$ws.Range("E92:G95").Font.Size = 11
$ws.Range("D92:D95").Font.Size = 12
$ws.Range("A92:G95").RowHeight = 18

# Row 92
$ws.Cells.Item(92, 1).Value = "Test Report "
$ws.Cells.Item(92, 2).Value = "V1"
$ws.Cells.Item(92, 3).Value = 82
$ws.Cells.Item(92, 4).Value = "testcases title is not descriptive enough "
$ws.Cells.Item(92, 5).Value = 45064
$ws.Cells.Item(92, 6).Value = "Tarek"
$ws.Cells.Item(92, 7).Value = "Follow Feature"

# Row 93
$ws.Cells.Item(93, 1).Value = "Test Report "
$ws.Cells.Item(93, 2).Value = "V1"
$ws.Cells.Item(93, 3).Value = 83
$ws.Cells.Item(93, 4).Value = 'add this step in the first testcase "1.Find the category page that the user wants to unfollow."'
$ws.Cells.Item(93, 5).Value = 45064
$ws.Cells.Item(93, 6).Value = "Tarek"
$ws.Cells.Item(93, 7).Value = "Follow Feature"

# Row 94
$ws.Cells.Item(94, 1).Value = "Test Report "
$ws.Cells.Item(94, 2).Value = "V1"
$ws.Cells.Item(94, 3).Value = 84
$ws.Cells.Item(94, 4).Value = "testcases title is not descriptive enough "
$ws.Cells.Item(94, 5).Value = 45064
$ws.Cells.Item(94, 6).Value = "Tarek"
$ws.Cells.Item(94, 7).Value = "notification Feature"

# Row 95
$ws.Cells.Item(95, 1).Value = "Test Report "
$ws.Cells.Item(95, 2).Value = "V1"
$ws.Cells.Item(95, 3).Value = 85
$ws.Cells.Item(95, 4).Value = "the steps of all testcases are not follow the same sequence"
$ws.Cells.Item(95, 5).Value = 45064
$ws.Cells.Item(95, 6).Value = "Tarek"
$ws.Cells.Item(95, 7).Value = "notification Feature"

# ------------------------------------------------------------------
# Final cursor position, matching where the author ended up after
# finishing the edits.
# ------------------------------------------------------------------
$ws.Range("G100").Select() | Out-Null
